# Finished the last open tasks:
# - Search Series and Users
# - Close Requests
# - Add Comments to Episodes
#
# Updates the "Tabelle1" sheet of WebViews.xlsx:
#  - Change Username/Password rows now reference the new user/[id]/* endpoints
#  - Mark several "Backend Completed" (F column) cells as done ("X") for the
#    rows that were finished: Add Episode Comment, Search Series, Change
#    Username/Password, Search Users, Elevate User, Close Request
#  - Move the active selection to C11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22/23: "Change Username" / "Change Password" now call the new
# user/[id]/name and user/[id]/password endpoints instead of the old
# currentUser/username and currentUser/password ones.
$ws.Range("C23").Value = "user/[id]/password"
$ws.Range("C22").Value = "user/[id]/name"

# Mark the now-finished rows as backend-completed ("X" in column F):
#  18 -> Add Episode Comment
#  20 -> Search Series
#  22 -> Change Username
#  23 -> Change Password
#  25 -> Search Users
#  27 -> Elevate User
#  31 -> Close Request
$ws.Range("F18").Value = "X"
$ws.Range("F20").Value = "X"
$ws.Range("F22").Value = "X"
$ws.Range("F23").Value = "X"
$ws.Range("F25").Value = "X"
$ws.Range("F27").Value = "X"
$ws.Range("F31").Value = "X"

# Move the selection to C11, matching where the author was last working.
$null = $ws.Range("C11").Select()
